# Separates tests from project framework
# The "SignIn" sheet had a "Url" column (column A) pointing to a hard-coded
# test environment URL (http://192.168.99.100:5000/). This column is removed
# so the test data no longer embeds an environment-specific URL, and the
# remaining Username/Password columns shift left.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("SignIn")
$ws.Activate()

$ws.Columns.Item(1).Delete()
